# This script applies the textual updates described by the diff:
# it updates the date line and the 25 "A×B=C" multiplication answers
# found throughout the document's table cells.

$d = $word.ActiveDocument

# Mapping of old text -> new text, in document order.
$replacements = [ordered]@{
    "2025-01-05 Sunday" = "2025-01-06 Monday"
    "536×3=1608"        = "857×7=5999"
    "365×2=730"         = "960×7=6720"
    "565×9=5085"        = "462×5=2310"
    "732×3=2196"        = "976×7=6832"
    "290×2=580"         = "987×8=7896"
    "937×5=4685"        = "224×3=672"
    "974×3=2922"        = "624×9=5616"
    "486×2=972"         = "881×8=7048"
    "841×9=7569"        = "743×3=2229"
    "428×8=3424"        = "700×3=2100"
    "295×4=1180"        = "114×4=456"
    "188×5=940"         = "431×4=1724"
    "136×7=952"         = "952×2=1904"
    "186×5=930"         = "387×3=1161"
    "778×2=1556"        = "349×4=1396"
    "856×3=2568"        = "795×5=3975"
    "107×7=749"         = "484×8=3872"
    "218×2=436"         = "667×9=6003"
    "308×6=1848"        = "428×9=3852"
    "110×9=990"         = "557×8=4456"
    "643×8=5144"        = "707×4=2828"
    "309×8=2472"        = "356×2=712"
    "334×3=1002"        = "475×3=1425"
    "601×9=5409"        = "969×7=6783"
    "676×5=3380"        = "893×6=5358"
}

foreach ($old in $replacements.Keys) {
    $new = $replacements[$old]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
